$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to rounded (2 decimal place) figures - "custom accuracy"
$ws.Range("B5").Value = 6.25
$ws.Range("C5").Value = 4.36
$ws.Range("D5").Value = 0.77
$ws.Range("E5").Value = 13.29
$ws.Range("F5").Value = 11.03
$ws.Range("G5").Value = 4.92
$ws.Range("H5").Value = 21.29
$ws.Range("I5").Value = 7.56
$ws.Range("J5").Value = 3.25
$ws.Range("K5").Value = 4.93
$ws.Range("L5").Value = 5.32
$ws.Range("M5").Value = 5.56
$ws.Range("N5").Value = 1.57
$ws.Range("O5").Value = 4.89
$ws.Range("P5").Value = 6.88
$ws.Range("Q5").Value = 4.26
$ws.Range("R5").Value = 0.75
$ws.Range("S5").Value = 0.41
$ws.Range("T5").Value = 67.07
$ws.Range("U5").Value = 13.76
$ws.Range("V5").Value = 4.51
$ws.Range("W5").Value = 9.08
$ws.Range("X5").Value = 4.93
$ws.Range("Y5").Value = 0.59
$ws.Range("Z5").Value = 9.86
$ws.Range("AA5").Value = 3.98
$ws.Range("AB5").Value = 3.64
$ws.Range("AC5").Value = 4.26
$ws.Range("AD5").Value = 5.59
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 19.26
$ws.Range("AG5").Value = 2.47
$ws.Range("AH5").Value = 5.64

# Remove row 6 (data trimmed to 1000 rows worth of source data -> fewer rows here)
$ws.Rows.Item(6).Delete()
